# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Note: several "Price" column values look numeric (e.g. "0.9996", "306.30")
# but must stay plain text, matching the exact digits/trailing zeros seen on
# coinranking.com - a leading apostrophe forces Excel to keep them as text
# instead of auto-converting to a Number and dropping significant trailing
# zeros / re-formatting them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.244.61'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '1.904.98'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('D4').Value = '''0.9996'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''306.30'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').Value = '''0.5423'
$ws.Range('E7').Value = '  +3.75%  '
$ws.Range('D8').Value = '''0.3809'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '''0.07299'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').Value = '''22.21'
$ws.Range('E10').Value = '  +5.31%  '
$ws.Range('D11').Value = '''0.9046'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '''0.08185'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '''95.87'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').Value = '''5.350'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '''0.9997'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '''0.000008651'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').Value = '''0.9995'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '27.265.13'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '1.175.55'
$ws.Range('E20').Value = '  -38.04%  '
$ws.Range('D21').Value = '''5.049'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = '''10.83'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = '''6.522'
$ws.Range('E23').Value = '  +1.93%  '
$ws.Range('B24').Value = 'LidoDAOToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D24').Value = '''2.312'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''148.52'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '''18.40'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '''1.756'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').Value = '''116.74'
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('D29').Value = '''4.866'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '''4.702'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').Value = '''0.09230'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '''0.8320'
$ws.Range('E32').Value = '  +5.54%  '
$ws.Range('D33').Value = '''0.05086'
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('D34').Value = '''1.223'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D35').Value = '''3.003'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('D36').Value = '''3.321'
$ws.Range('E36').Value = '  -2.91%  '
$ws.Range('D37').Value = '''2.698'
$ws.Range('E37').Value = '  +4.42%  '
$ws.Range('D38').Value = '''0.5945'
$ws.Range('E38').Value = '  +4.43%  '
$ws.Range('D39').Value = '''0.02002'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').Value = '''1.081'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '''9.294'
$ws.Range('E41').Value = '  +3.04%  '
$ws.Range('D42').Value = '''6.679'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').Value = '''116.58'
$ws.Range('D44').Value = '''0.5129'
$ws.Range('E44').Value = '  +5.67%  '
$ws.Range('D45').Value = '''0.1531'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('D47').Value = '''0.9989'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('D48').Value = '''1.644'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D49').Value = '''38.32'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = '''0.06103'
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('D51').Value = '''63.60'
$ws.Range('E51').Value = '  +0.28%  '
